$wb = $excel.ActiveWorkbook

# --- Sheet "Present-Storage" ---
$ws1 = $wb.Worksheets.Item("Present-Storage")

# Column B width change (9.6 -> 8.4). The interop engine snaps ColumnWidth
# to the nearest 1/6-character increment (pixel grid), so 7.5 is the input
# that lands closest to the target stored width of 8.4.
$ws1.Range("B:B").ColumnWidth = 7.5

$ws1.Range("B2").Value = 7.71
$ws1.Range("B3").Value = 5.53
$ws1.Range("B4").Value = 8.33
$ws1.Range("B5").Value = 8.24
$ws1.Range("B6").Value = 7
$ws1.Range("B7").Value = 4.7
$ws1.Range("B8").Value = 7.66
$ws1.Range("B9").Value = 7.56
$ws1.Range("B10").Value = 15.8
$ws1.Range("B11").Value = 14.11
$ws1.Range("B12").Value = 16.29
$ws1.Range("B13").Value = 16.22
$ws1.Range("B14").Value = 7.95
$ws1.Range("B15").Value = 5.86
$ws1.Range("B16").Value = 8.539999999999999
$ws1.Range("B17").Value = 8.449999999999999

# --- Sheet "2030-Storage" ---
$ws2 = $wb.Worksheets.Item("2030-Storage")

$ws2.Range("B3").Value = 4.41
$ws2.Range("B7").Value = 4.14
$ws2.Range("B11").Value = 8.41
$ws2.Range("B15").Value = 5.55

# --- Sheet "2050-Storage" ---
$ws3 = $wb.Worksheets.Item("2050-Storage")

$ws3.Range("B3").Value = 68.48999999999999
$ws3.Range("B7").Value = 71.90000000000001
$ws3.Range("B11").Value = 58.41
$ws3.Range("B15").Value = 67.03
